# "fixed export and fixing maps"
#
# The published table had an extra subtitle row (census-results caveat) and
# two extra year columns (1989 / 2002) that were dropped from the export,
# leaving only the 2014 figures. The worksheet tab is also renamed from the
# generic "1" to the municipality name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab: "1" -> "ბაღდათი"
$ws.Name = "ბაღდათი"

# Remove the "(მოსახლეობის აღწერის შედეგებით)" subtitle row entirely -
# everything below shifts up one row.
$ws.Rows("2").Delete()

# Keep only the 2014 figures - drop the 1989 and 2002 columns (B:C),
# leaving the old column D (2014) as the new column B.
$ws.Columns("B:C").Delete()

# Match the saved selection cursor.
$ws.Range("A2").Select() | Out-Null
